$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:U) to (B:V)
$ws.Columns("A:A").Insert()

# Add the new "SlNo" header and sequential numbers for the data rows
$ws.Range("A1").Value = "SlNo"
$ws.Range("A2").Value = 10000
$ws.Range("A3").Value = 10001
$ws.Range("A4").Value = 10002
$ws.Range("A5").Value = 10003

# Update the active selection to reflect where the user ended up after editing
$ws.Range("A6").Select()
